{"js": "// Split the run \"do C\u1ee5c CSQLHC v\u1ec1 TTXH c\u1ea5p\" into three runs:\n//   \"do \" + \"C\u1ee5c CSQLHC v\u1ec1 TTXH\" + \" c\u1ea5p\"\n// All three keep the exact same run formatting (color 000000) as the\n// original single run \u2014 only the run boundaries change so that the\n// middle phrase \"C\u1ee5c CSQLHC v\u1ec1 TTXH\" becomes its own run.\n\nconst body = context.document.body;\n\n// Locate the middle phrase inside the target sentence. The surrounding\n// text (\"do \" / \" c\u1ea5p\") is unique enough in this document that searching\n// directly for the phrase is safe, but scope the search to the exact\n// sentence first to be extra precise.\nconst sentence = body.search(\"do C\u1ee5c CSQLHC v\u1ec1 TTXH c\u1ea5p\", { matchCase: true });\nsentence.load(\"items\");\nawait context.sync();\n\nif (sentence.items.length === 0) {\n  throw new Error(\"Target sentence not found\");\n}\n\nconst target = sentence.items[0];\nconst middle = target.search(\"C\u1ee5c CSQLHC v\u1ec1 TTXH\", { matchCase: true });\nmiddle.load(\"items\");\nawait context.sync();\n\nif (middle.items.length === 0) {\n  throw new Error(\"Target phrase not found inside sentence\");\n}\n\nconst middleRange = middle.items[0];\n\n// Toggling a character property on the sub-range and then flipping it\n// back to its original value forces Word to materialize the sub-range\n// as its own run (splitting the original run into three), while the\n// resulting formatting stays identical to the source run.\nmiddleRange.font.bold = true;\nawait context.sync();\n\nmiddleRange.font.bold = false;\nawait context.sync();\n", "ps1": "# Split the run \"do C\u1ee5c CSQLHC v\u1ec1 TTXH c\u1ea5p\" into three runs:\n#   \"do \" + \"C\u1ee5c CSQLHC v\u1ec1 TTXH\" + \" c\u1ea5p\"\n# All three keep the exact same run formatting (color 000000) as the\n# original single run -- only the run boundaries change so that the\n# middle phrase \"C\u1ee5c CSQLHC v\u1ec1 TTXH\" becomes its own run.\n\n$d = $word.ActiveDocument\n\n# Locate the full sentence first (it is unique in the document), then\n# locate the middle phrase inside that sentence so we don't accidentally\n# match text elsewhere.\n$sentence = $d.Content\n$foundSentence = $sentence.Find.Execute(\"do C\u1ee5c CSQLHC v\u1ec1 TTXH c\u1ea5p\")\nif (-not $foundSentence) {\n    throw \"Target sentence 'do C\u1ee5c CSQLHC v\u1ec1 TTXH c\u1ea5p' not found\"\n}\n$sentStart = $sentence.Start\n$sentEnd = $sentence.End\n\n$middle = $d.Range($sentStart, $sentEnd)\n$foundMiddle = $middle.Find.Execute(\"C\u1ee5c CSQLHC v\u1ec1 TTXH\")\nif (-not $foundMiddle) {\n    throw \"Target phrase 'C\u1ee5c CSQLHC v\u1ec1 TTXH' not found inside sentence\"\n}\n$midStart = $middle.Start\n$midEnd = $middle.End\n\n# Toggling a character property on the sub-range and then flipping it back\n# to its original value forces Word to materialize the sub-range as its\n# own run (splitting the original run into three), while the resulting\n# formatting stays identical to the source run.\n$d.Range($midStart, $midEnd).Font.Bold = 1\n$d.Range($midStart, $midEnd).Font.Bold = 0\n"}
